$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Check register")

# --- Row height adjustments -------------------------------------------------
# Header spacer row gets taller
$ws.Rows.Item(1).RowHeight = 34
# Thin separator row just above the table header grows by 1pt
$ws.Rows.Item(9).RowHeight = 11
# The whole transactions table (header row 10 through the last filled data
# row 33) grows to accommodate additional wrapped text in the new rows
$ws.Range("A10:A33").EntireRow.RowHeight = 65

# --- New test-case rows (33 & 34) ------------------------------------------
# Row 33: new wishlist test case
$ws.Range("C33").Value = 45800
$ws.Range("D33").Value = "1001"
$ws.Range("G33").Value = "위시리스트 하트 버튼 눌러 추가하기"
$ws.Range("K33").Value = "failed"

# Row 34: new color-selection test case
$ws.Range("C34").Value = 45800
$ws.Range("D34").Value = "1002"
$ws.Range("G34").Value = "사용자 선택 색상 테스트"

# --- View state --------------------------------------------------------------
# Scroll the visible window down / left and move the selection to A24
$window = $excel.ActiveWindow
$window.ScrollRow = 18
$window.ScrollColumn = 1
$ws.Range("A24").Select() | Out-Null
